$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove underline from the two hyperlink-styled cells (D2, F2) while
# keeping the hyperlink colour - touched individually so no stray empty
# cell gets materialised in between (E2).
$ws.Range("D2").Font.Underline = $False
$ws.Range("F2").Font.Underline = $False

# Widen column D slightly
$ws.Columns("D").ColumnWidth = 18.33

# Add a new test case row: TC_0056fcwf / Chrome, matching the formatting
# already used by the rows above it.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "TC_0056fcwf"
$ws.Range("B7").Value = "Chrome"

# Move the active selection to B13
$ws.Range("B13").Select()
